$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the latest "Censos" (census) data point, 2023,
# right after the existing census series (old row 7 shifts down to 8).
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "Censos"
$ws.Range("B7").Value = 2023
$ws.Range("C7").Value = 3444.2629999999999

# Append the new "Proy (NNUU 2019)" (UN projection) data point for 2020
# at the bottom of the table.
$ws.Range("A15").Value = "Proy (NNUU 2019) "
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = 3473.7269999999999

# Column A best-fits to the longest label now in the sheet.
$ws.Columns.Item(1).AutoFit()

# Update the current selection to match the edited file.
$ws.Range("C8").Select()
